$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.229.92"
$ws.Range("E2").Value = "  +5.47%  "

$ws.Range("D3").Value = "'3.611.04"
$ws.Range("E3").Value = "  +5.22%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "'591.49"
$ws.Range("E5").Value = "  +3.88%  "

$ws.Range("D6").Value = "'189.80"
$ws.Range("E6").Value = "  +3.66%  "

$ws.Range("D7").Value = "'0.645"
$ws.Range("E7").Value = "  +2.24%  "

$ws.Range("D8").Value = "'3.605.30"
$ws.Range("E8").Value = "  +5.25%  "

$ws.Range("E9").Value = "  +0.04%  "

$ws.Range("E10").Value = "  +3.20%  "

$ws.Range("D11").Value = "'0.664"
$ws.Range("E11").Value = "  +3.47%  "

$ws.Range("D12").Value = "'58.36"
$ws.Range("E12").Value = "  +4.15%  "

$ws.Range("D13").Value = "'0.0000289"
$ws.Range("E13").Value = "  +3.90%  "

$ws.Range("D14").Value = "'9.87"
$ws.Range("E14").Value = "  +5.21%  "

$ws.Range("D15").Value = "'4.190.44"
$ws.Range("E15").Value = "  +5.27%  "

$ws.Range("D16").Value = "'19.66"
$ws.Range("E16").Value = "  +6.03%  "

$ws.Range("D17").Value = "'3.613.10"
$ws.Range("E17").Value = "  +5.32%  "

$ws.Range("D18").Value = "'70.240.18"
$ws.Range("E18").Value = "  +5.57%  "

$ws.Range("D19").Value = "'12.60"
$ws.Range("E19").Value = "  +4.58%  "

$ws.Range("E20").Value = "  +0.66%  "

$ws.Range("E21").Value = "  +4.21%  "

$ws.Range("D22").Value = "'488.30"
$ws.Range("E22").Value = "  +0.43%  "

$ws.Range("D23").Value = "'18.60"
$ws.Range("E23").Value = "  +13.44%  "

$ws.Range("E24").Value = "  +7.21%  "

$ws.Range("E25").Value = "  +6.39%  "

$ws.Range("D26").Value = "'90.81"
$ws.Range("E26").Value = "  +2.08%  "

$ws.Range("E27").Value = "  +5.67%  "

$ws.Range("D28").Value = "'11.16"
$ws.Range("E28").Value = "  +2.12%  "

$ws.Range("D29").Value = "'9.44"
$ws.Range("E29").Value = "  +4.29%  "

$ws.Range("D30").Value = "'32.83"
$ws.Range("E30").Value = "  +4.91%  "

$ws.Range("D31").Value = "'7.76"
$ws.Range("E31").Value = "  +8.16%  "

$ws.Range("D32").Value = "'12.34"
$ws.Range("E32").Value = "  +5.56%  "

$ws.Range("D33").Value = "'624.82"
$ws.Range("E33").Value = "  +5.30%  "

$ws.Range("D34").Value = "'0.119"
$ws.Range("E34").Value = "  +6.66%  "

$ws.Range("D35").Value = "'65.63"
$ws.Range("E35").Value = "  +3.83%  "

$ws.Range("D36").Value = "'0.0₃0819"
$ws.Range("E36").Value = "  +6.85%  "

$ws.Range("D37").Value = "'38.48"
$ws.Range("E37").Value = "  +6.47%  "

$ws.Range("B38").Value = "Dai"
$ws.Range("C38").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = "  +0.06%  "

$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "'0.147"
$ws.Range("E39").Value = "  -0.02%  "

$ws.Range("E40").Value = "  +4.25%  "

$ws.Range("D41").Value = "'3.59"
$ws.Range("E41").Value = "  +0.77%  "

$ws.Range("D42").Value = "'3.313.46"
$ws.Range("E42").Value = "  +5.23%  "

$ws.Range("D43").Value = "'3.15"
$ws.Range("E43").Value = "  +7.66%  "

$ws.Range("E44").Value = "  +6.12%  "

$ws.Range("D45").Value = "'2.70"
$ws.Range("E45").Value = "  +6.38%  "

$ws.Range("E46").Value = "  +2.68%  "

$ws.Range("E47").Value = "  +3.49%  "

$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D48").Value = "'9.07"
$ws.Range("E48").Value = "  +4.07%  "

$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D49").Value = "'2.72"
$ws.Range("E49").Value = "  -1.66%  "

$ws.Range("D50").Value = "'3.31"
$ws.Range("E50").Value = "  +5.40%  "

$ws.Range("B51").Value = "FirstDigitalUSD"
$ws.Range("C51").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D51").Value = "'0.999"
$ws.Range("E51").Value = "  +0.04%  "
